$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing the existing rows 14.. down by one
# (old row 14 becomes 15, ..., old row 60 becomes the new row 61).
$ws.Rows(14).Insert()

# The new row 14 is a fresh weekly record. Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T
# are constant across the whole data block, so copy them from the row right
# below (the former row 14, now shifted to row 15).
$ws.Range("A14:T14").Value = $ws.Range("A15:T15").Value()

# Now fill in the values that differ for this new record.
$ws.Range("D14").Value = 44690   # Fecha
$ws.Range("M14").Value = 35      # Volumen
$ws.Range("N14").Value = 30000   # Precio minimo
$ws.Range("O14").Value = 30000   # Precio maximo
$ws.Range("P14").Value = 30000   # Precio promedio ponderado
$ws.Range("S14").Value = 1500    # Precio $/Kg
